$d = $word.ActiveDocument

$old = "Read through your forms and make sure everything is correct. "
$new = "Read through your forms and make sure everything is correct. You do not need to fill out the " + [char]0x201C + "Affix Transfer Tax Stamp" + [char]0x201D + " box until later."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
